# "creating properties Type: SAVE."
#
# A new "properties" row (properties / origin / Deviation) is inserted
# above the Fica-rate table on the "Configs" sheet. Inserting the row
# pushes the whole table (including the merged Year column B4:B18)
# down by one row, exactly like choosing "Insert Sheet Rows" in Excel
# on row 3. A new, empty trailing row is then added at the bottom of
# the (now one-row-taller) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# --- 1. Insert a new row above the current header row -----------------
# Rows 3:18 (data + the merged Year groups) shift down to rows 4:19;
# mergeCells (B4:B6, B7:B9, ...) move down automatically with them.
$ws.Rows("3:3").Insert()

# --- 2. Populate the new "properties" row ------------------------------
$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"
$ws.Range("E3").Value = ""

$ws.Range("D3").Borders.Item(10).LineStyle = 1
$ws.Range("D3").Borders.Item(10).Weight = 2

# --- 3. Add a new, blank trailing row under the table ------------------
$lastRow = 20
$ws.Range("B" + $lastRow).Value = ""
$ws.Range("C" + $lastRow).Value = ""
$ws.Range("D" + $lastRow).Value = ""
$ws.Range("E" + $lastRow).Value = ""

$ws.Range("B" + $lastRow).Borders.Item(7).LineStyle = 1
$ws.Range("B" + $lastRow).Borders.Item(7).Weight = 2
$ws.Range("D" + $lastRow).Borders.Item(10).LineStyle = 1
$ws.Range("D" + $lastRow).Borders.Item(10).Weight = 2
